$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target state for the data table (rows 16-22), columns B..G
# B = Tipo Doc Trabajador, C = N Doc Trabajador, D = Nombre Trabajador,
# E = Periodo Mora, F = Salario Basico, G = Valor Mora

# Row 16: MIRLLAN YULIETH MARTINEZ PITALUA - periodo 1906
$ws.Cells.Item(16, 2).Value = "CC"
$ws.Cells.Item(16, 3).Value = "1047446850"
$ws.Cells.Item(16, 4).Value = "MIRLLAN YULIETH MARTINEZ PITALUA"
$ws.Cells.Item(16, 5).Value = "1906"
$ws.Cells.Item(16, 6).Value = 33125
$ws.Cells.Item(16, 7).Value = 1423500

# Row 17: GLADYS DEL SOCORRO LOPEZ LLERENA - periodo 1906
$ws.Cells.Item(17, 2).Value = "CC"
$ws.Cells.Item(17, 3).Value = "45442352"
$ws.Cells.Item(17, 4).Value = "GLADYS DEL SOCORRO LOPEZ LLERENA"
$ws.Cells.Item(17, 5).Value = "1906"
$ws.Cells.Item(17, 6).Value = 33125
$ws.Cells.Item(17, 7).Value = 828116

# Row 18: MIRLLAN YULIETH MARTINEZ PITALUA - periodo 1907
$ws.Cells.Item(18, 2).Value = "CC"
$ws.Cells.Item(18, 3).Value = "1047446850"
$ws.Cells.Item(18, 4).Value = "MIRLLAN YULIETH MARTINEZ PITALUA"
$ws.Cells.Item(18, 5).Value = "1907"
$ws.Cells.Item(18, 6).Value = 33125
$ws.Cells.Item(18, 7).Value = 1423500

# Row 19: GLADYS DEL SOCORRO LOPEZ LLERENA - periodo 1907
$ws.Cells.Item(19, 2).Value = "CC"
$ws.Cells.Item(19, 3).Value = "45442352"
$ws.Cells.Item(19, 4).Value = "GLADYS DEL SOCORRO LOPEZ LLERENA"
$ws.Cells.Item(19, 5).Value = "1907"
$ws.Cells.Item(19, 6).Value = 33125
$ws.Cells.Item(19, 7).Value = 828116

# Row 20: MIRLLAN YULIETH MARTINEZ PITALUA - periodo 1908
$ws.Cells.Item(20, 2).Value = "CC"
$ws.Cells.Item(20, 3).Value = "1047446850"
$ws.Cells.Item(20, 4).Value = "MIRLLAN YULIETH MARTINEZ PITALUA"
$ws.Cells.Item(20, 5).Value = "1908"
$ws.Cells.Item(20, 6).Value = 33125
$ws.Cells.Item(20, 7).Value = 1423500

# Row 21: GLADYS DEL SOCORRO LOPEZ LLERENA - periodo 1908
$ws.Cells.Item(21, 2).Value = "CC"
$ws.Cells.Item(21, 3).Value = "45442352"
$ws.Cells.Item(21, 4).Value = "GLADYS DEL SOCORRO LOPEZ LLERENA"
$ws.Cells.Item(21, 5).Value = "1908"
$ws.Cells.Item(21, 6).Value = 8833
$ws.Cells.Item(21, 7).Value = 828116

# Row 22: MIRLLAN YULIETH MARTINEZ PITALUA - periodo 1909
$ws.Cells.Item(22, 2).Value = "CC"
$ws.Cells.Item(22, 3).Value = "1047446850"
$ws.Cells.Item(22, 4).Value = "MIRLLAN YULIETH MARTINEZ PITALUA"
$ws.Cells.Item(22, 5).Value = "1909"
$ws.Cells.Item(22, 6).Value = 33125
$ws.Cells.Item(22, 7).Value = 1423500

$wb.Save()
